# Remove the trailing "Ver no Jupiter..." and "© 2020 ..." paragraphs
# (and the blank paragraph immediately preceding them) that follow the
# bibliography entry ending in "CRC Press USA:1997. 3040P.".
#
# Document tail structure (before):
#   ... CRC Press USA:1997. 3040P.
#   <blank paragraph>
#   Ver no Jupiter Salvar em pdf Salvar em docx
#   © 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution
#   <blank paragraph>
#   <page-break paragraph>
#
# Desired structure (after):
#   ... CRC Press USA:1997. 3040P.
#   <blank paragraph>
#   <page-break paragraph>

$d = $word.ActiveDocument

$startIndex = $null
$endIndex = $null
$index = 0

foreach ($p in $d.Paragraphs) {
    $index = $index + 1
    $t = $p.Range.Text
    if ($t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startIndex = $index
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $endIndex = $index
    }
}

if ($startIndex -ne $null -and $endIndex -ne $null) {
    # Also remove the blank paragraph immediately preceding the "Ver no
    # Jupiter" paragraph so we don't leave a duplicate blank line behind.
    if ($startIndex -gt 1) {
        $precedingPara = $d.Paragraphs.Item($startIndex - 1)
        if ($precedingPara.Range.Text -eq "`r") {
            $startIndex = $startIndex - 1
        }
    }

    $rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End

    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
